$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 corresponds to "Biceps curl" - add description, sets, reps
$ws.Range("B4").Value = "Begin by standing up holding a dumbell in each of your hands with your arms handing by your sides. Keep your palms facing forward during the exercise. Keep your upper arm(s) in place as you curl the dumbells up to your shoulders."
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 10

# Row 7 corresponds to "Deadlift" - add description, sets, reps
$ws.Range("B7").Value = "Begin by positioning the barbell on the floor in front of you and loading the desired amount of weight. Stand facing the barbell with your feet shoulder width apart, and your toes under the barbell. Keeping your feet flat, down and grab the barbell with your hands at a shoulder-width distance apart. Lift the barbell while keeping it close to your legs, with your shoulders back and your chest up. Try to prevent your back from rounding. Lift the bar to thigh level then return to the ground, in the same position it started."
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 5

# The new wrapped descriptions grow those two rows - match Excel's own autosize
$ws.Rows("4").RowHeight = 75
$ws.Rows("7").RowHeight = 174.75

# Update selection/view state to match target
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B11").Select()
